# Weekly update: a new daily observation is inserted as a new row 76,
# pushing the existing rows 76-171 down to 77-172.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 76 (shifts rows 76:171 -> 77:172)
$ws.Rows(76).Insert()

# Fill in the new row 76 with the new data point
$ws.Cells.Item(76,1).Value = 7
$ws.Cells.Item(76,2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(76,3).Value = 'Ñuble'
$ws.Cells.Item(76,4).Value = 44789
$ws.Cells.Item(76,5).Value = 16
$ws.Cells.Item(76,6).Value = 100112045
$ws.Cells.Item(76,7).Value = 'Zapallo'
$ws.Cells.Item(76,8).Value = 'Camote'
$ws.Cells.Item(76,9).Value = '1a (guarda)'
$ws.Cells.Item(76,10).Value = 200
$ws.Cells.Item(76,11).Value = 850
$ws.Cells.Item(76,12).Value = 900
$ws.Cells.Item(76,13).Value = 875
$ws.Cells.Item(76,14).Value = '$/kilo (volumen en unidades)'
$ws.Cells.Item(76,15).Value = "Región de O'Higgins"
$ws.Cells.Item(76,16).Value = 875
$ws.Cells.Item(76,17).Value = 1
$ws.Cells.Item(76,18).Value = 'Hortaliza'
